$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Manually restore the AGN identifier for two epochs that had previously
# been treated as spurious and left blank in column A (they belong to
# the same AGN as the row immediately above, just a different epoch).
$ws.Range("A6").Value2 = $ws.Range("A5").Value2
$ws.Range("A12").Value2 = $ws.Range("A11").Value2

# Row 12 was incorrectly tagged as the W3 band (2nd max/min band, which
# is no longer being used) - correct it to W2.
$ws.Range("B12").Value2 = "W2"

# Leave the active selection where the user finished working.
$ws.Range("F18").Select() | Out-Null
